$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 16:35"

# Row 4 - Estados Unidos: refreshed case counts
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 1623352
$ws.Range("C4").Value = 2450
$ws.Range("D4").Value = 382944
$ws.Range("E4").Value = 1143976
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 78
$ws.Range("H4").Value = 96432

# Row 55 - Noruega: refreshed case counts
$ws.Range("A55").Value = "Noruega"
$ws.Range("B55").Value = 8309
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 7727
$ws.Range("E55").Value = 347
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 235

# Row 89 - Cuba: refreshed case counts
$ws.Range("A89").Value = "Cuba"
$ws.Range("B89").Value = 1916
$ws.Range("C89").Value = 8
$ws.Range("D89").Value = 1631
$ws.Range("E89").Value = 204
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 81

# Sierra Leona is now added with updated stats, pushing Malta and
# Republica del Chad down one row each (table is sorted by total cases).
# Row 126 - Sierra Leona (new position/data)
$ws.Range("A126").Value = "Sierra Leona"
$ws.Range("B126").Value = 606
$ws.Range("C126").Value = 21
$ws.Range("D126").Value = 230
$ws.Range("E126").Value = 338
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 3
$ws.Range("H126").Value = 38

# Row 127 - Malta (shifted down, data unchanged)
$ws.Range("A127").Value = "Malta"
$ws.Range("B127").Value = 600
$ws.Range("C127").Value = 1
$ws.Range("D127").Value = 469
$ws.Range("E127").Value = 125
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 6

# Row 128 - Republica del Chad (shifted down, data unchanged)
$ws.Range("A128").Value = "Republica del Chad"
$ws.Range("B128").Value = 588
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 186
$ws.Range("E128").Value = 344
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 58

# Groenlandia is now added with updated stats, pushing Montserrat and
# Seychelles down one row each.
# Row 209 - Groenlandia (new position/data)
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("B209").Value = 11
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 11
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

# Row 210 - Montserrat (shifted down, data unchanged)
$ws.Range("A210").Value = "Montserrat"
$ws.Range("B210").Value = 11
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 10
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 1

# Row 211 - Seychelles (shifted down, data unchanged)
$ws.Range("A211").Value = "Seychelles"
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 11
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0
